$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing PATIENT_A header (column B) to PATIENT_A_SnpEff,
# then insert a new column C for PATIENT_A_dbNSFP (shifting the old
# SnpEff_overall_impact_rank column from C to D).
$ws.Range("B1").Value = "PATIENT_A_SnpEff"

$ws.Columns("C").Insert()

$ws.Range("C1").Value = "PATIENT_A_dbNSFP"
$ws.Range("C2").Value = 12
$ws.Range("C4").Value = 4

# New cell style for C2: larger font (12pt) in the existing blue
# (#6699FF) used elsewhere, on a new orange (#FFA500) fill.
$ws.Range("C2").Font.Size = 12
$ws.Range("C2").Font.Color = 0xFF9966
$ws.Range("C2").Interior.Color = 0x00A5FF

# New cell style for C4: larger font (12pt) in the existing dark blue
# (#003366) used elsewhere, on a new (white) fill.
$ws.Range("C4").Font.Size = 12
$ws.Range("C4").Font.Color = 0x663300
$ws.Range("C4").Interior.Color = 0xFFFFFF
